# Remove the "lemmalist-greek" dependency row from the Acknowledgments sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Acknowledgments")

# Find the row whose "Name" column (A) contains "lemmalist-greek" and delete it entirely.
$found = $ws.Cells.Find("lemmalist-greek", [Type]::Missing, [Type]::Missing, 1)
if ($found -ne $null) {
    $found.EntireRow.Delete()
}
